$wb = $excel.ActiveWorkbook

# @@ -3730,25 +3730,25 @@  sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5061618.5
$ws.Range("J62").Value = 15311.25
$ws.Range("L62").Value = 15311.25
$ws.Range("N62").Value = -16559.25

# @@ -3883,25 +3883,25 @@  sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5061618.5
$ws.Range("J65").Value = 15311.25
$ws.Range("L65").Value = 76556.25
$ws.Range("N65").Value = -82796.25

# @@ -10725,22 +10725,22 @@  sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3208
$ws.Range("I61").Value = 2523.3684
$ws.Range("K61").Value = 2523.3684
$ws.Range("M61").Value = -2311.3684

# @@ -11353,22 +11353,22 @@  sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5474.893
$ws.Range("I74").Value = 1011.88
$ws.Range("K74").Value = 1011.88
$ws.Range("M74").Value = -137.88

# @@ -11500,22 +11500,22 @@  sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5474.893
$ws.Range("I77").Value = 1011.88
$ws.Range("K77").Value = 5059.4
$ws.Range("M77").Value = -691.3999999999996

# @@ -14183,25 +14183,25 @@  sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2895.8445
$ws.Range("I132").Value = 2417.8108
$ws.Range("J132").Value = 5106.75
$ws.Range("K132").Value = 7253.432400000001
$ws.Range("L132").Value = 15320.25
$ws.Range("M132").Value = -4723.432400000001
$ws.Range("N132").Value = -20380.25

# @@ -14379,22 +14379,22 @@  sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3208
$ws.Range("I136").Value = 2523.3684
$ws.Range("K136").Value = 7570.1052
$ws.Range("M136").Value = -5020.1052

# @@ -19248,25 +19248,25 @@  sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 917.3333
$ws.Range("I94").Value = 789.913
$ws.Range("J94").Value = 1650
$ws.Range("K94").Value = 789.913
$ws.Range("L94").Value = 1650
$ws.Range("M94").Value = -338.913
$ws.Range("N94").Value = -2552

# @@ -26172,25 +26172,25 @@  sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2616
$ws.Range("I94").Value = 1650
$ws.Range("J94").Value = 2892
$ws.Range("K94").Value = 1650
$ws.Range("L94").Value = 2892
$ws.Range("M94").Value = -1199
$ws.Range("N94").Value = -3794

# @@ -28141,25 +28141,25 @@  sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2517.2703
$ws.Range("I134").Value = 1411.25
$ws.Range("J134").Value = 5958.222
$ws.Range("K134").Value = 4233.75
$ws.Range("L134").Value = 17874.666
$ws.Range("M134").Value = -1698.75
$ws.Range("N134").Value = -22944.666

# @@ -28627,25 +28627,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 90909310
$ws.Range("I2").Value = 19.5
$ws.Range("J2").Value = 142857490
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 857144940
$ws.Range("M2").Value = -4
$ws.Range("N2").Value = -857145166

# @@ -28783,25 +28783,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1188.05
$ws.Range("I5").Value = 521.2727
$ws.Range("J5").Value = 2003
$ws.Range("K5").Value = 1563.8181
$ws.Range("L5").Value = 6009
$ws.Range("M5").Value = -1451.8181
$ws.Range("N5").Value = -6233

# @@ -29646,25 +29646,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 900
$ws.Range("J22").Value = 1800
$ws.Range("L22").Value = 5400
$ws.Range("N22").Value = -5738

# @@ -29897,25 +29897,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 900
$ws.Range("J27").Value = 1800
$ws.Range("L27").Value = 5400
$ws.Range("N27").Value = -5604

# @@ -30249,25 +30249,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1698.2
$ws.Range("I34").Value = 1245.5
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 3736.5
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -3652.5
$ws.Range("N34").Value = -6168

# @@ -30503,22 +30503,22 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9170
$ws.Range("J39").Value = 9170
$ws.Range("L39").Value = 27510
$ws.Range("N39").Value = -28098

# @@ -31308,22 +31308,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1950
$ws.Range("I55").Value = 900
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 2700
$ws.Range("L55").Value = 9000
$ws.Range("M55").Value = -2523
$ws.Range("N55").Value = -9354

# @@ -31709,22 +31712,22 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4405.5
$ws.Range("I63").Value = 4405.5
$ws.Range("K63").Value = 13216.5
$ws.Range("M63").Value = -12467.5

# @@ -31862,22 +31865,22 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 4405.5
$ws.Range("I66").Value = 4405.5
$ws.Range("K66").Value = 39649.5
$ws.Range("M66").Value = -35905.5

# @@ -34277,25 +34280,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 877.25
$ws.Range("I114").Value = 1028
$ws.Range("J114").Value = 847.1
$ws.Range("K114").Value = 3084
$ws.Range("L114").Value = 2541.3
$ws.Range("M114").Value = 170
$ws.Range("N114").Value = -9049.299999999999

# @@ -35143,25 +35146,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5209638.5
$ws.Range("J131").Value = 5377673.5
$ws.Range("L131").Value = 16133020.5
$ws.Range("N131").Value = -16143100.5

# @@ -35351,25 +35354,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1188.05
$ws.Range("I135").Value = 521.2727
$ws.Range("J135").Value = 2003
$ws.Range("K135").Value = 4691.454299999999
$ws.Range("L135").Value = 18027
$ws.Range("M135").Value = -2156.454299999999
$ws.Range("N135").Value = -23097

# @@ -35455,25 +35458,25 @@  sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 10104253
$ws.Range("I137").Value = 20003912
$ws.Range("J137").Value = 204593.2
$ws.Range("K137").Value = 60011736
$ws.Range("L137").Value = 613779.6000000001
$ws.Range("M137").Value = -60006636
$ws.Range("N137").Value = -623979.6000000001

# @@ -42137,25 +42140,25 @@  sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3467.625
$ws.Range("I132").Value = 2909.4119
$ws.Range("J132").Value = 4823.2856
$ws.Range("K132").Value = 8728.235700000001
$ws.Range("L132").Value = 14469.8568
$ws.Range("M132").Value = -6198.235700000001
$ws.Range("N132").Value = -19529.8568

# @@ -47195,22 +47198,25 @@  sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2874.0557
$ws.Range("I93").Value = 2966.6667
$ws.Range("J93").Value = 2688.8333
$ws.Range("K93").Value = 2966.6667
$ws.Range("L93").Value = 2688.8333
$ws.Range("M93").Value = -1718.6667
$ws.Range("N93").Value = -5184.8333

# @@ -49079,25 +49085,25 @@  sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3776.647
$ws.Range("I132").Value = 2241.2
$ws.Range("J132").Value = 5970.143
$ws.Range("K132").Value = 6723.599999999999
$ws.Range("L132").Value = 17910.429
$ws.Range("M132").Value = -4193.599999999999
$ws.Range("N132").Value = -22970.429

# @@ -49131,22 +49137,22 @@  sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 49979.816
$ws.Range("J133").Value = 49979.816
$ws.Range("L133").Value = 49979.816
$ws.Range("N133").Value = -55039.816

# @@ -49278,25 +49284,25 @@  sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5284.7617
$ws.Range("I136").Value = 2842.1936
$ws.Range("J136").Value = 12168.363
$ws.Range("K136").Value = 8526.5808
$ws.Range("L136").Value = 36505.089
$ws.Range("M136").Value = -5976.5808
$ws.Range("N136").Value = -41605.089

# @@ -53585,25 +53591,25 @@  sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3935.4722
$ws.Range("I81").Value = 1755.909
$ws.Range("J81").Value = 4894.48
$ws.Range("K81").Value = 3511.818
$ws.Range("L81").Value = 9788.959999999999
$ws.Range("M81").Value = -2450.818
$ws.Range("N81").Value = -11910.96

# @@ -53729,25 +53735,25 @@  sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3935.4722
$ws.Range("I84").Value = 1755.909
$ws.Range("J84").Value = 4894.48
$ws.Range("K84").Value = 17559.09
$ws.Range("L84").Value = 48944.8
$ws.Range("M84").Value = -12255.09
$ws.Range("N84").Value = -59552.8

# @@ -54305,25 +54311,25 @@  sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1739.6
$ws.Range("I96").Value = 924.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 924.5
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = 448.5
$ws.Range("N96").Value = -7746

# @@ -56051,25 +56057,25 @@  sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2370.6099
$ws.Range("I132").Value = 2060.3076
$ws.Range("J132").Value = 2908.4666
$ws.Range("K132").Value = 6180.9228
$ws.Range("L132").Value = 8725.399800000001
$ws.Range("M132").Value = -3650.9228
$ws.Range("N132").Value = -13785.3998
